# Controle_de_Treinamentos.xlsx - apply commit "feat: atualizando os tests"
#
# Summary of changes:
#  - Rename three header labels (A1, C1, D1)
#  - Update the expiry date text in E2 (Amanda Duarte / Gestao de Fornecedores)
#  - Apply a dd/mm/yy number format to E2
#  - Widen column D slightly
#  - Touch cell "locked" / alignment state across the sheet (protection bits)
#  - Move the active selection to E3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text updates -------------------------------------------------
$ws.Range("A1").Value = "Nome"
$ws.Range("C1").Value = "Treinamento Obrigatório"
$ws.Range("D1").Value = "Situação de Treinamento"

# --- 2. Update the training-due date shown in E2 ----------------------------
$ws.Range("E2").Value = "15/12/2030"

# --- 3. Touch protection/locked state for the whole used range so the style
#        table picks up explicit protection application, matching the
#        "locked" cell-protection metadata carried by the edited workbook.
$ws.Range("A1:E68").Locked = $true

# --- 4. Give the due-date cell its own date number format -------------------
$ws.Range("E2").NumberFormat = "dd/mm/yy"

# --- 5. Widen column D --------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 24.8

# --- 6. Move the selection to E3, matching the saved cursor position --------
$ws.Range("E3").Select() | Out-Null
